$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (preserve exact string, e.g. trailing
# zeros / multi-dot "thousand" separators) without leaving the cell tagged
# with a quote-prefix style -- matches the original inlineStr formatting.
function Set-TextValue($range, $text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

# --- Simple value updates (price / volume columns) ---
Set-TextValue $ws.Range("D2") "62.801.77"
$ws.Range("E2").Value = "  -1.85%  "
Set-TextValue $ws.Range("D3") "2.567.45"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  +0.04%  "
Set-TextValue $ws.Range("D5") "567.20"
$ws.Range("E5").Value = "  -1.67%  "
Set-TextValue $ws.Range("D6") "153.09"
$ws.Range("E6").Value = "  -2.64%  "
Set-TextValue $ws.Range("D7") "0.999"
$ws.Range("E7").Value = "  +0.01%  "
Set-TextValue $ws.Range("D8") "0.617"
$ws.Range("E8").Value = "  -1.09%  "
Set-TextValue $ws.Range("D9") "2.569.89"
$ws.Range("E9").Value = "  -1.68%  "
$ws.Range("E10").Value = "  -5.12%  "
Set-TextValue $ws.Range("D11") "5.62"
$ws.Range("E11").Value = "  -3.53%  "
$ws.Range("E12").Value = "  -0.10%  "
Set-TextValue $ws.Range("D13") "0.374"
$ws.Range("E13").Value = "  -2.59%  "
Set-TextValue $ws.Range("D14") "27.58"
$ws.Range("E14").Value = "  -2.71%  "
Set-TextValue $ws.Range("D15") "3.029.00"
$ws.Range("E15").Value = "  -1.99%  "
Set-TextValue $ws.Range("D16") "0.0000177"
Set-TextValue $ws.Range("D17") "62.681.15"
$ws.Range("E17").Value = "  -1.59%  "
Set-TextValue $ws.Range("D18") "2.553.86"
$ws.Range("E18").Value = "  -2.07%  "
Set-TextValue $ws.Range("D19") "11.78"
$ws.Range("E19").Value = "  -2.52%  "
Set-TextValue $ws.Range("D20") "7.34"
$ws.Range("E20").Value = "  -4.02%  "
Set-TextValue $ws.Range("D21") "4.40"
$ws.Range("E21").Value = "  -3.94%  "
Set-TextValue $ws.Range("D22") "333.50"
$ws.Range("E22").Value = "  -2.95%  "
$ws.Range("E23").Value = "  +0.02%  "
Set-TextValue $ws.Range("D24") "66.97"
$ws.Range("E24").Value = "  -0.50%  "
$ws.Range("E25").Value = "  +5.20%  "
$ws.Range("E26").Value = "  -2.18%  "
Set-TextValue $ws.Range("D27") "2.703.56"
$ws.Range("E27").Value = "  -1.59%  "
Set-TextValue $ws.Range("D28") "8.95"
$ws.Range("E28").Value = "  -2.96%  "
$ws.Range("E29").Value = "  -0.14%  "
Set-TextValue $ws.Range("D30") "554.92"
$ws.Range("E30").Value = "  -7.23%  "
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("E34").Value = "  -4.25%  "
Set-TextValue $ws.Range("D35") "1.67"
$ws.Range("E35").Value = "  -4.63%  "
Set-TextValue $ws.Range("D36") "6.36"
$ws.Range("E36").Value = "  -3.76%  "
Set-TextValue $ws.Range("D37") "5.16"
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("E38").Value = "  -0.06%  "
Set-TextValue $ws.Range("D39") "0.394"
$ws.Range("E39").Value = "  -3.49%  "
Set-TextValue $ws.Range("D40") "19.30"
$ws.Range("E40").Value = "  -2.58%  "
Set-TextValue $ws.Range("D41") "152.04"
$ws.Range("E41").Value = "  -1.42%  "
Set-TextValue $ws.Range("D42") "1.83"
$ws.Range("E42").Value = "  -2.78%  "
Set-TextValue $ws.Range("D43") "0.999"
$ws.Range("E43").Value = "  -0.02%  "
Set-TextValue $ws.Range("D44") "2.43"
$ws.Range("E44").Value = "  -1.67%  "
Set-TextValue $ws.Range("D45") "155.04"
$ws.Range("E45").Value = "  -0.68%  "
Set-TextValue $ws.Range("D46") "23.15"
$ws.Range("E46").Value = "  +0.23%  "
Set-TextValue $ws.Range("D47") "3.81"
$ws.Range("E47").Value = "  -3.19%  "
Set-TextValue $ws.Range("D48") "0.0570"
$ws.Range("E48").Value = "  -3.95%  "
Set-TextValue $ws.Range("D49") "0.619"
$ws.Range("E49").Value = "  -1.69%  "
Set-TextValue $ws.Range("D50") "0.0981"
$ws.Range("E50").Value = "  -3.74%  "
Set-TextValue $ws.Range("D51") "0.0244"
$ws.Range("E51").Value = "  -1.89%  "

# --- Row 31/32 swap: Aptos now ranks above Binance-PegBSC-USD ---
$ws.Range("B31").Value = "Aptos"
$ws.Range("C31").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D31") "7.96"
$ws.Range("E31").Value = "  +0.49%  "
$ws.Range("B32").Value = "Binance-PegBSC-USD"
$ws.Range("C32").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D32") "1.00"
$ws.Range("E32").Value = "  -0.01%  "
